$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct writes. D-column price values are forced to
# text via a leading apostrophe (classic Excel "treat as text" marker)
# so purely-numeric-looking strings (e.g. "1.00", "0.0000279") are not
# silently coerced into numbers by Excel's type inference. The style
# is then reset to "Normal" to drop the implicit text-number-format
# that the apostrophe trick applies, keeping cell styling identical
# to the original (unstyled) cells.

$ws.Range("D2").Formula = "'65.491.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").Formula = "'3.390.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Formula = "'559.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Formula = "'175.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("D8").Formula = "'3.379.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Formula = "'0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.21%  "

$ws.Range("D11").Formula = "'0.638"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("D12").Formula = "'53.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("D13").Formula = "'0.0000279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").Formula = "'9.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.58%  "

$ws.Range("D15").Formula = "'3.933.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Formula = "'3.378.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").Formula = "'65.282.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").Formula = "'11.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("D21").Formula = "'1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("D22").Formula = "'468.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").Formula = "'4.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.30%  "

$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Formula = "'14.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.18%  "

$ws.Range("D26").Formula = "'87.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.37%  "

$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").Formula = "'10.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").Formula = "'8.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").Formula = "'31.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.14%  "

$ws.Range("D31").Formula = "'6.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.92%  "

$ws.Range("D32").Formula = "'63.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.89%  "

$ws.Range("D33").Formula = "'11.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").Formula = "'573.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Formula = "'3.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("D38").Formula = "'0.141"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Formula = "'35.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("D40").Formula = "'0.374"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.02%  "

$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").Formula = "'3.121.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("D45").Formula = "'3.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("E47").Value = "  -2.56%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").Formula = "'140.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "

$ws.Range("D50").Formula = "'2.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("D51").Formula = "'8.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
